$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9999999435414485
$ws.Range("E2").Value = 0.9999999435414485

# Row 3
$ws.Range("D3").Value = 0.9999999987657784
$ws.Range("E3").Value = 0.9999999987657784

# Row 4
$ws.Range("D4").Value = 0.000000001345302800396394
$ws.Range("E4").Value = 0.000000001345302800396394

# Row 5
$ws.Range("D5").Value = 0.9993654803362613
$ws.Range("E5").Value = 0.9993654803362613

# Row 6 (Success flips FALSE -> TRUE)
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = 0.04485462019112854
$ws.Range("E6").Value = 0.04485462019112854

# Row 7
$ws.Range("D7").Value = 0.9999999976328175
$ws.Range("E7").Value = 0.000000002367182450058181

# Row 8
$ws.Range("D8").Value = 0.9998470511737478
$ws.Range("E8").Value = 0.0001529488262521639

# Row 9
$ws.Range("D9").Value = 0.000000103776861125329
$ws.Range("E9").Value = 0.9999998962231389

# Row 10
$ws.Range("D10").Value = 0.999805785079903
$ws.Range("E10").Value = 0.0001942149200969689

# Row 11
$ws.Range("D11").Value = 0.0009588548240338831
$ws.Range("E11").Value = 0.9990411451759661
$ws.Range("F11").Value = 6.764225959777832
$ws.Range("G11").Value = 0.5

